# NGFS Corridor Maps Sketch Data - input data update
# Adds a new arterial ("Tasman Dr + Montague Expy + Comstock St + McCarthy
# Blvd + S Main St") replacing the old "Tasman Dr + Montague Expy" text for
# the AM SouthBay_Inner rows, plus three new link rows on the
# arterial_link_pick sheet and four corresponding rows on the am_links
# lookup sheet.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("arterial_link_pick")
$ws3 = $wb.Worksheets.Item("am_links")

# --- arterial_link_pick: rename the arterial text for the existing AM rows
$newArterial = "Tasman Dr + Montague Expy + Comstock St + McCarthy Blvd + S Main St"
$ws2.Range("I13").Value = $newArterial
$ws2.Range("I14").Value = $newArterial
$ws2.Range("I15").Value = $newArterial
$ws2.Range("I16").Value = $newArterial

# --- am_links: append four new lookup rows (46-49)
$ws3.Range("A46").Value = "5328_5312"
$ws3.Range("B46").Value = "SouthBay_Inner"

$ws3.Range("A47").Value = "5709_5674"
$ws3.Range("B47").Value = "SouthBay_Inner"

$ws3.Range("A48").Value = "5801_4369"
$ws3.Range("B48").Value = "SouthBay_Inner"

$ws3.Range("A49").Value = "5799_5101"
$ws3.Range("C49").Value = "SouthBay_Inner"

# --- arterial_link_pick: append three new rows (55-57)
$ws2.Range("A55").Value = "Santa Clara"
$ws2.Range("B55").Value = 6
$ws2.Range("C55").Value = "SouthBay_Inner"
$ws2.Range("D55").Value = "5602_5653"
$ws2.Range("E55").Value = 5602
$ws2.Range("F55").Value = 5653
$ws2.Range("G55").Value = "AM"
$ws2.Range("H55").Value = 2
$ws2.Range("I55").Value = $newArterial
$ws2.Range("J55").Value = 5328
$ws2.Range("K55").Value = 5312

$ws2.Range("A56").Value = "Santa Clara"
$ws2.Range("B56").Value = 6
$ws2.Range("C56").Value = "SouthBay_Inner"
$ws2.Range("D56").Value = "5799_5101"
$ws2.Range("E56").Value = 5799
$ws2.Range("F56").Value = 5101
$ws2.Range("G56").Value = "AM"
$ws2.Range("H56").Value = 3
$ws2.Range("I56").Value = $newArterial
$ws2.Range("J56").Value = 5709
$ws2.Range("K56").Value = 5674

$ws2.Range("A57").Value = "Santa Clara"
$ws2.Range("B57").Value = 6
$ws2.Range("C57").Value = "SouthBay_Inner"
$ws2.Range("D57").Value = "5799_5101"
$ws2.Range("E57").Value = 5799
$ws2.Range("F57").Value = 5101
$ws2.Range("G57").Value = "AM"
$ws2.Range("H57").Value = 3
$ws2.Range("I57").Value = $newArterial
$ws2.Range("J57").Value = 5801
$ws2.Range("K57").Value = 4369

# fill the L column (a_b concat) formula down for the new rows
$ws2.Range("L55:L57").Formula = '=_xlfn.CONCAT(J55,"_",K55)'

# --- restore the view/selection state (best effort)
$ws2.Range("H58").Select()
$ws3.Range("B51").Select()
$ws3.Select()
